# Weekly update: insert the newest price observation for this market/product
# (Hortaliza, Vega Monumental Concepción - Ají) as a new row 123, pushing the
# existing rows 123:187 down to 124:188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 123 (shifts rows 123..187 down to 124..188,
# and extends the used range from A1:R187 to A1:R188).
$ws.Rows(123).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A123").Value = 11
$ws.Range("B123").Value = "Vega Monumental Concepción"
$ws.Range("C123").Value = "Bíobío"
$ws.Range("D123").Value = 45029
$ws.Range("E123").Value = 8
$ws.Range("F123").Value = 100112021
$ws.Range("G123").Value = "Ají"
$ws.Range("H123").Value = "Americana (o)"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 50
$ws.Range("K123").Value = 17000
$ws.Range("L123").Value = 18000
$ws.Range("M123").Value = 17400
$ws.Range("N123").Value = "$/saco 25 kilos"
$ws.Range("O123").Value = "Región Metropolitana"
$ws.Range("P123").Value = 696
$ws.Range("Q123").Value = 25
$ws.Range("R123").Value = "Hortaliza"
